$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 30.75612566666667
$ws.Range("H2").Value = 92.268377
$ws.Range("I2").Value = 0.9777985798685588
$ws.Range("J2").Value = 0.9777985798685588
$ws.Range("M2").Value = 49.89274333333334
$ws.Range("N2").Value = 149.67823
$ws.Range("O2").Value = 0.8663408689480834
$ws.Range("P2").Value = 0.8663408689480835
$ws.Range("Q2").Value = 1534.507483814746
$ws.Range("R2").Value = 13810.56735433271
$ws.Range("S2").Value = 0.8471068713395291
$ws.Range("T2").Value = 0.8471068713395292
# Row 3
$ws.Range("G3").Value = 30.75612566666667
$ws.Range("H3").Value = 92.268377
$ws.Range("I3").Value = 0.9777985798685588
$ws.Range("J3").Value = 0.9777985798685588
$ws.Range("O3").Value = 0.06984725491313053
$ws.Range("P3").Value = 0.06984725491313053
$ws.Range("Q3").Value = 123.7170486003449
$ws.Range("R3").Value = 1113.453437403104
$ws.Range("S3").Value = 0.06829654666177624
$ws.Range("T3").Value = 0.06829654666177624
# Row 4
$ws.Range("G4").Value = 30.75612566666667
$ws.Range("H4").Value = 92.268377
$ws.Range("I4").Value = 0.9777985798685588
$ws.Range("J4").Value = 0.9777985798685588
$ws.Range("M4").Value = 1.266267666666667
$ws.Range("N4").Value = 3.798803
$ws.Range("O4").Value = 0.02198755485004457
$ws.Range("P4").Value = 0.02198755485004457
$ws.Range("Q4").Value = 38.94548748363678
$ws.Range("R4").Value = 350.509387352731
$ws.Range("S4").Value = 0.02149939990715562
$ws.Range("T4").Value = 0.02149939990715562
# Row 5
$ws.Range("G5").Value = 30.75612566666667
$ws.Range("H5").Value = 92.268377
$ws.Range("I5").Value = 0.9777985798685588
$ws.Range("J5").Value = 0.9777985798685588
$ws.Range("M5").Value = 0.2206823333333333
$ws.Range("N5").Value = 0.6620469999999999
$ws.Range("O5").Value = 0.003831942516052412
$ws.Range("P5").Value = 0.003831942516052413
$ws.Range("Q5").Value = 6.787333576413221
$ws.Range("R5").Value = 61.086002187719
$ws.Range("S5").Value = 0.003746867950334001
$ws.Range("T5").Value = 0.003746867950334001
# Row 6
$ws.Range("G6").Value = 30.75612566666667
$ws.Range("H6").Value = 92.268377
$ws.Range("I6").Value = 0.9777985798685588
$ws.Range("J6").Value = 0.9777985798685588
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 2.187988666666667
$ws.Range("N6").Value = 6.563966
$ws.Range("O6").Value = 0.03799237877268909
$ws.Range("P6").Value = 0.03799237877268909
$ws.Range("Q6").Value = 67.29405438924243
$ws.Range("R6").Value = 605.646489503182
$ws.Range("S6").Value = 0.03714889400976377
$ws.Range("T6").Value = 0.03714889400976377
# Row 7
$ws.Range("I7").Value = 0.004830327290741966
$ws.Range("J7").Value = 0.004830327290741966
$ws.Range("M7").Value = 49.89274333333334
$ws.Range("N7").Value = 149.67823
$ws.Range("O7").Value = 0.8663408689480834
$ws.Range("P7").Value = 0.8663408689480835
$ws.Range("Q7").Value = 7.580470589264445
$ws.Range("R7").Value = 68.22423530338001
$ws.Range("S7").Value = 0.004184709942365036
$ws.Range("T7").Value = 0.004184709942365037
# Row 8
$ws.Range("I8").Value = 0.004830327290741966
$ws.Range("J8").Value = 0.004830327290741966
$ws.Range("O8").Value = 0.06984725491313053
$ws.Range("P8").Value = 0.06984725491313053
$ws.Range("S8").Value = 0.0003373851015903052
$ws.Range("T8").Value = 0.0003373851015903052
# Row 9
$ws.Range("I9").Value = 0.004830327290741966
$ws.Range("J9").Value = 0.004830327290741966
$ws.Range("M9").Value = 1.266267666666667
$ws.Range("N9").Value = 3.798803
$ws.Range("O9").Value = 0.02198755485004457
$ws.Range("P9").Value = 0.02198755485004457
$ws.Range("Q9").Value = 0.1923908000242222
$ws.Range("R9").Value = 1.731517200218
$ws.Range("S9").Value = 0.0001062070862488561
$ws.Range("T9").Value = 0.0001062070862488562
# Row 10
$ws.Range("I10").Value = 0.004830327290741966
$ws.Range("J10").Value = 0.004830327290741966
$ws.Range("M10").Value = 0.2206823333333333
$ws.Range("N10").Value = 0.6620469999999999
$ws.Range("O10").Value = 0.003831942516052412
$ws.Range("P10").Value = 0.003831942516052413
$ws.Range("Q10").Value = 0.03352944387577778
$ws.Range("R10").Value = 0.301764994882
$ws.Range("S10").Value = 0.0000185095365118424
$ws.Range("T10").Value = 0.0000185095365118424
# Row 11
$ws.Range("I11").Value = 0.004830327290741966
$ws.Range("J11").Value = 0.004830327290741966
$ws.Range("K11").Value = 3.0
$ws.Range("L11").Value = 1.0
$ws.Range("M11").Value = 2.187988666666667
$ws.Range("N11").Value = 6.563966
$ws.Range("O11").Value = 0.03799237877268909
$ws.Range("P11").Value = 0.03799237877268909
$ws.Range("Q11").Value = 0.3324327873995556
$ws.Range("R11").Value = 2.991895086596
$ws.Range("S11").Value = 0.0001835156240259259
$ws.Range("T11").Value = 0.0001835156240259259
# Row 12
$ws.Range("E12").Value = 3.0
$ws.Range("F12").Value = 1.0
$ws.Range("G12").Value = 0.4673666666666667
$ws.Range("H12").Value = 1.4021
$ws.Range("I12").Value = 0.01485851852399773
$ws.Range("J12").Value = 0.01485851852399773
$ws.Range("M12").Value = 49.89274333333334
$ws.Range("N12").Value = 149.67823
$ws.Range("O12").Value = 0.8663408689480834
$ws.Range("P12").Value = 0.8663408689480835
$ws.Range("Q12").Value = 23.31820514255556
$ws.Range("R12").Value = 209.863846283
$ws.Range("S12").Value = 0.01287254184936139
$ws.Range("T12").Value = 0.01287254184936139
# Row 13
$ws.Range("E13").Value = 3.0
$ws.Range("F13").Value = 1.0
$ws.Range("G13").Value = 0.4673666666666667
$ws.Range("H13").Value = 1.4021
$ws.Range("I13").Value = 0.01485851852399773
$ws.Range("J13").Value = 0.01485851852399773
$ws.Range("O13").Value = 0.06984725491313053
$ws.Range("P13").Value = 0.06984725491313053
$ws.Range("Q13").Value = 1.879990517688889
$ws.Range("R13").Value = 16.9199146592
$ws.Range("S13").Value = 0.001037826730977142
$ws.Range("T13").Value = 0.001037826730977142
# Row 14
$ws.Range("E14").Value = 3.0
$ws.Range("F14").Value = 1.0
$ws.Range("G14").Value = 0.4673666666666667
$ws.Range("H14").Value = 1.4021
$ws.Range("I14").Value = 0.01485851852399773
$ws.Range("J14").Value = 0.01485851852399773
$ws.Range("M14").Value = 1.266267666666667
$ws.Range("N14").Value = 3.798803
$ws.Range("O14").Value = 0.02198755485004457
$ws.Range("P14").Value = 0.02198755485004457
$ws.Range("Q14").Value = 0.5918112984777779
$ws.Range("R14").Value = 5.326301686300001
$ws.Range("S14").Value = 0.0003267024910368034
$ws.Range("T14").Value = 0.0003267024910368035
# Row 15
$ws.Range("E15").Value = 3.0
$ws.Range("F15").Value = 1.0
$ws.Range("G15").Value = 0.4673666666666667
$ws.Range("H15").Value = 1.4021
$ws.Range("I15").Value = 0.01485851852399773
$ws.Range("J15").Value = 0.01485851852399773
$ws.Range("M15").Value = 0.2206823333333333
$ws.Range("N15").Value = 0.6620469999999999
$ws.Range("O15").Value = 0.003831942516052412
$ws.Range("P15").Value = 0.003831942516052413
$ws.Range("Q15").Value = 0.1031395665222222
$ws.Range("R15").Value = 0.9282560987
$ws.Range("S15").Value = 0.00005693698885765925
$ws.Range("T15").Value = 0.00005693698885765925
# Row 16
$ws.Range("E16").Value = 3.0
$ws.Range("F16").Value = 1.0
$ws.Range("G16").Value = 0.4673666666666667
$ws.Range("H16").Value = 1.4021
$ws.Range("I16").Value = 0.01485851852399773
$ws.Range("J16").Value = 0.01485851852399773
$ws.Range("K16").Value = 3.0
$ws.Range("L16").Value = 1.0
$ws.Range("M16").Value = 2.187988666666667
$ws.Range("N16").Value = 6.563966
$ws.Range("O16").Value = 0.03799237877268909
$ws.Range("P16").Value = 0.03799237877268909
$ws.Range("Q16").Value = 1.022592969844444
$ws.Range("R16").Value = 9.2033367286
$ws.Range("S16").Value = 0.000564510463764739
$ws.Range("T16").Value = 0.000564510463764739
# Row 17
$ws.Range("E17").Value = 2.0
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.07903166666666667
$ws.Range("H17").Value = 0.237095
$ws.Range("I17").Value = 0.002512574316701549
$ws.Range("J17").Value = 0.002512574316701549
$ws.Range("M17").Value = 49.89274333333334
$ws.Range("N17").Value = 149.67823
$ws.Range("O17").Value = 0.8663408689480834
$ws.Range("P17").Value = 0.8663408689480835
$ws.Range("Q17").Value = 3.943106660205556
$ws.Range("R17").Value = 35.48795994185
$ws.Range("S17").Value = 0.002176745816827857
$ws.Range("T17").Value = 0.002176745816827857
# Row 18
$ws.Range("E18").Value = 2.0
$ws.Range("F18").Value = 0.6666666666666666
$ws.Range("G18").Value = 0.07903166666666667
$ws.Range("H18").Value = 0.237095
$ws.Range("I18").Value = 0.002512574316701549
$ws.Range("J18").Value = 0.002512574316701549
$ws.Range("O18").Value = 0.06984725491313053
$ws.Range("P18").Value = 0.06984725491313053
$ws.Range("Q18").Value = 0.3179062490488889
$ws.Range("R18").Value = 2.86115624144
$ws.Range("S18").Value = 0.0001754964187868379
$ws.Range("T18").Value = 0.0001754964187868379
# Row 19
$ws.Range("E19").Value = 2.0
$ws.Range("F19").Value = 0.6666666666666666
$ws.Range("G19").Value = 0.07903166666666667
$ws.Range("H19").Value = 0.237095
$ws.Range("I19").Value = 0.002512574316701549
$ws.Range("J19").Value = 0.002512574316701549
$ws.Range("M19").Value = 1.266267666666667
$ws.Range("N19").Value = 3.798803
$ws.Range("O19").Value = 0.02198755485004457
$ws.Range("P19").Value = 0.02198755485004457
$ws.Range("Q19").Value = 0.1000752441427778
$ws.Range("R19").Value = 0.9006771972850001
$ws.Range("S19").Value = 0.00005524536560328856
$ws.Range("T19").Value = 0.00005524536560328857
# Row 20
$ws.Range("E20").Value = 2.0
$ws.Range("F20").Value = 0.6666666666666666
$ws.Range("G20").Value = 0.07903166666666667
$ws.Range("H20").Value = 0.237095
$ws.Range("I20").Value = 0.002512574316701549
$ws.Range("J20").Value = 0.002512574316701549
$ws.Range("M20").Value = 0.2206823333333333
$ws.Range("N20").Value = 0.6620469999999999
$ws.Range("O20").Value = 0.003831942516052412
$ws.Range("P20").Value = 0.003831942516052413
$ws.Range("Q20").Value = 0.01744089260722222
$ws.Range("R20").Value = 0.156968033465
$ws.Range("S20").Value = 0.000009628040348910004
$ws.Range("T20").Value = 0.000009628040348910006
# Row 21
$ws.Range("E21").Value = 2.0
$ws.Range("F21").Value = 0.6666666666666666
$ws.Range("G21").Value = 0.07903166666666667
$ws.Range("H21").Value = 0.237095
$ws.Range("I21").Value = 0.002512574316701549
$ws.Range("J21").Value = 0.002512574316701549
$ws.Range("K21").Value = 3.0
$ws.Range("L21").Value = 1.0
$ws.Range("M21").Value = 2.187988666666667
$ws.Range("N21").Value = 6.563966
$ws.Range("O21").Value = 0.03799237877268909
$ws.Range("P21").Value = 0.03799237877268909
$ws.Range("Q21").Value = 0.1729203909744444
$ws.Range("R21").Value = 1.55628351877
$ws.Range("S21").Value = 0.00009545867513465572
$ws.Range("T21").Value = 0.00009545867513465572
